$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the model values below
# (a routine holdings-data refresh) can be written, then restore protection
# with the same password.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclosure footer text (A9).
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."
$ws.Range("A9").Value = $newText

# Refresh the weight/percent-change model data (rows 2-6).
$ws.Range("D2").Value = 0.2483243590373695
$ws.Range("E2").Value = 0.01966085033177678

$ws.Range("D3").Value = 0.2484538030734802
$ws.Range("E3").Value = 0.02886771300448432

$ws.Range("D4").Value = 0.2578573148779001
$ws.Range("E4").Value = 0.0158029053788773

$ws.Range("D5").Value = 0.2453645230112503
$ws.Range("E5").Value = 0.02125603864734305

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.02134493367084689

$ws.Protect("D382")
